$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 with refreshed timestamp value (tiny precision change from the scheduled task run)
$ws.Range("A2").Value = 45865.04190887731

# Append new row 3 with the latest reading
$ws.Range("A3").Value = 45865.08354313159
$ws.Range("B3").Value = 2025
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 13.28
$ws.Range("E3").Value = 90.41
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 4.14
$ws.Range("H3").Value = "E"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "02:00:18"

# Match the date/time style used by column A on row 2
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
